$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.292.64'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = '2.583.12'
$ws.Range('E3').Value = '  -2.21%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '589.23'
$ws.Range('E5').Value = '  -3.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '150.59'
$ws.Range('E6').Value = '  +2.26%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '0.586'
$ws.Range('E8').Value = '  -0.54%  '
$ws.Range('E9').Value = '  +0.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '5.70'
$ws.Range('E10').Value = '  +1.78%  '
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range('D13').Value = '27.54'
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('D14').Value = '3.050.19'
$ws.Range('E14').Value = '  -2.06%  '
$ws.Range('D15').Value = '63.109.16'
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '0.0000155'
$ws.Range('E16').Value = '  +5.34%  '
$ws.Range('D17').Value = '2.558.17'
$ws.Range('E17').Value = '  -2.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '12.22'
$ws.Range('E18').Value = '  +3.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '4.73'
$ws.Range('E19').Value = '  +3.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '344.98'
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '6.86'
$ws.Range('E21').Value = '  -0.94%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '67.37'
$ws.Range('E23').Value = '  +1.62%  '
$ws.Range('E24').Value = '  +0.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '9.24'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '1.66'
$ws.Range('E26').Value = '  -2.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '559.86'
$ws.Range('E27').Value = '  -0.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '8.06'
$ws.Range('E28').Value = '  -0.76%  '
$ws.Range('E29').Value = '  +0.85%  '
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('D32').Value = '0.0₃0855'
$ws.Range('E32').Value = '  +0.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '1.75'
$ws.Range('E33').Value = '  -0.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '5.24'
$ws.Range('E34').Value = '  -0.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '166.66'
$ws.Range('E35').Value = '  -1.71%  '
$ws.Range('E36').Value = '  +2.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '19.52'
$ws.Range('E38').Value = '  +1.72%  '
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '166.55'
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '39.60'
$ws.Range('E42').Value = '  -1.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '3.99'
$ws.Range('E43').Value = '  +5.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '22.86'
$ws.Range('E44').Value = '  +4.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '0.0583'
$ws.Range('E45').Value = '  +2.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '2.11'
$ws.Range('E46').Value = '  +5.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '0.628'
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '0.0251'
$ws.Range('E48').Value = '  +2.18%  '
$ws.Range('E49').Value = '  +0.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '19.22'
$ws.Range('E50').Value = '  +2.04%  '
$ws.Range('D51').Value = '0.0₆0232'
$ws.Range('E51').Value = '  +17.66%  '
